$p = $ppt.ActivePresentation
$s = $p.Slides.Item(50)
$notes = $s.NotesPage
$shp = $notes.Shapes.Item(2)
$shp.TextFrame.TextRange.Text = ""
